# Added some notes and parts spreadsheets
# On the "interfaces" sheet, insert a new row before the existing
# "RESET_N" row (row 60) documenting the SRCLK2 / Load clock signal,
# shifting everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("interfaces")

$ws.Rows.Item(60).Insert()

$ws.Range("B60").Value = "SRCLK2"
$ws.Range("E60").Value = 1
$ws.Range("G60").Value = "Load clock for ser/des"

# Match the author's final selection/scroll position on the sheet.
$ws.Activate()
$null = $ws.Range("B61").Select()
